$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function New-Package([string]$bodyXml) {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
           '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
           '<pkg:xmlData>' + `
           '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" ' + `
           'xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' + `
           '<w:body>' + $bodyXml + '</w:body></w:document>' + `
           '</pkg:xmlData></pkg:part></pkg:package>'
}

function Find-ParagraphIndex([string]$startsWith) {
    $idx = 0
    foreach ($p in $d.Paragraphs) {
        $idx = $idx + 1
        if ($p.Range.Text.StartsWith($startsWith)) {
            return $idx
        }
    }
    return -1
}

# --- 1) Participants paragraph: split "Ceridwen" into its own run, ---
# --- wrapped in proofErr spellStart/spellEnd markers.               ---
$participantsIdx = Find-ParagraphIndex("Participants:")
$participantsRange = $d.Paragraphs($participantsIdx).Range

$participantsInner = `
    "<w:r $wNs><w:t xml:space=`"preserve`">Participants: </w:t></w:r>" + `
    "<w:r $wNs><w:tab/><w:t xml:space=`"preserve`">Shamim Bavani, Joanna Oruba, Samuel Coyle, </w:t></w:r>" + `
    "<w:proofErr $wNs w:type=`"spellStart`"/>" + `
    "<w:r $wNs><w:t>Ceridwen</w:t></w:r>" + `
    "<w:proofErr $wNs w:type=`"spellEnd`"/>" + `
    "<w:r $wNs><w:t xml:space=`"preserve`"> Grey, Joshua MacKay</w:t></w:r>"

$participantsPara = '<w:p w14:paraId="3E1651D4" w14:textId="7BD1064F" w:rsidR="00541E08" w:rsidRDefault="00541E08">' + `
    $participantsInner + '</w:p>'

[void]$participantsRange.InsertXML((New-Package $participantsPara), $null)

# --- 2) "Next meeting" paragraph: replace with the newly planned text ---
$meetingIdx = Find-ParagraphIndex("Next meeting")
$meetingRange = $d.Paragraphs($meetingIdx).Range

$meetingInner = `
    "<w:r $wNs><w:t>Next</w:t></w:r>" + `
    "<w:r $wNs><w:t xml:space=`"preserve`"> meeting</w:t></w:r>" + `
    "<w:r $wNs><w:t xml:space=`"preserve`"> planned:</w:t></w:r>" + `
    "<w:r $wNs><w:t xml:space=`"preserve`"> straight after interview or</w:t></w:r>" + `
    "<w:r $wNs><w:t xml:space=`"preserve`"> teams </w:t></w:r>" + `
    "<w:r $wNs><w:t>on</w:t></w:r>" + `
    "<w:r $wNs><w:t xml:space=`"preserve`"> Thursday</w:t></w:r>" + `
    "<w:r $wNs><w:t xml:space=`"preserve`">, </w:t></w:r>" + `
    "<w:r $wNs><w:t>16/02/2023.</w:t></w:r>" + `
    "<w:r $wNs><w:t xml:space=`"preserve`"> </w:t></w:r>"

$meetingPara = '<w:p w14:paraId="19372923" w14:textId="7B0A6E55" w:rsidR="00A448F1" w:rsidRDefault="00A448F1" w:rsidP="00A448F1">' + `
    $meetingInner + '</w:p>'

[void]$meetingRange.InsertXML((New-Package $meetingPara), $null)
